# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect refreshed scrape counts.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F (same rows/values on both sheets)
$updates = @{
    2  = 320
    4  = 10374
    6  = 945
    8  = 7183
    10 = 440
    11 = 204
    13 = 3200
    15 = 314
    16 = 661
    18 = 1036
    20 = 69
    21 = 1639
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
